$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.453128576278687
$ws.Range("B1").Value = 5.465088367462158
$ws.Range("C1").Value = 2.52459979057312
$ws.Range("D1").Value = 1.693475484848022
$ws.Range("E1").Value = 1.662680387496948
